$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.444.16"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.829.03"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'314.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.5158"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.04%  "
$ws.Range("D8").Value = "'0.3921"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.07669"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'41.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").Value = "'1.111"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").Value = "'21.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("D13").Value = "'6.283"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "'7.557"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "1.826.79"
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "'93.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.25%  "
$ws.Range("D18").Value = "'0.00001093"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").Value = "'0.06725"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").Value = "'17.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'6.194"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").Value = "28.458.46"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "'11.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "'2.252"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.93%  "
$ws.Range("D26").Value = "'156.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "'20.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "2.037.92"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").Value = "'2.402"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").Value = "'124.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "'0.1089"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "'5.676"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").Value = "'3.663"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'0.07007"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "'0.2223"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'8.933"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.49%  "
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").Value = "'5.147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "'0.6286"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("D41").Value = "'11.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "'1.183"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").Value = "'13.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "'0.5900"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").Value = "'3.709"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "'124.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'1.977"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "'1.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("D51").Value = "'0.06937"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.98%  "
